$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Hedged in townhouse" -> "Hedged" + "-" + "in townhouse"
#    (three runs sharing the same run formatting, space replaced by a
#    hyphen). Locate the exact standalone paragraph (NotesToBeDeleted
#    style, text is exactly "Hedged in townhouse") so the other, unrelated
#    occurrence of similar text elsewhere in the document is left alone.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Hedged in townhouse`r") {
        $start = $p.Range.Start

        # Replace just the single space right after "Hedged" with "-".
        $spaceRng = $d.Range($start + 6, $start + 7)
        $spaceRng.Text = "-"

        # Nudge a character-formatting property on just the new "-" run
        # (and revert it) so the run is kept distinct from its neighbours
        # instead of being silently coalesced back with them.
        $dashRng = $d.Range($start + 6, $start + 7)
        $dashRng.Bold = 1
        $dashRng2 = $d.Range($start + 6, $start + 7)
        $dashRng2.Bold = 0

        break
    }
}

# ---------------------------------------------------------------------
# 2) "<Color ,mood, Signposting, Narrative >" note paragraph: collapse
#    the three runs (separated by gramStart/gramEnd proofing marks) into
#    a single run with no proofing marks.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Color ,mood*") {
        $rng = $p.Range
        $rng.MoveEnd(1, -1)
        $full = $rng.Text
        $rng.Text = ""
        $insertPoint = $d.Range($rng.Start, $rng.Start)
        $insertPoint.InsertAfter($full)
        break
    }
}

# ---------------------------------------------------------------------
# 3) Footer DATE field cached text: 4/21/2021 -> 4/27/2021
# ---------------------------------------------------------------------
foreach ($sec in $d.Sections) {
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            $ftr.Range.Find.Execute("4/21/2021", $false, $false, $false, $false, $false, `
                $true, 1, $false, "4/27/2021", 2) | Out-Null
        }
    }
}
